# Agregar talleres 3 y 4
# (commit adds a "Taller 5" block: TIRM theory stub + a worked TIRM example
#  with a small cash-flow table and an NPV() formula)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Taller 5 header: "Tasa Interna de Retorno TIRM" (row 40) ---
$ws.Range("A40").Value = 5
$ws.Range("A40").Font.Bold = $true
$ws.Range("C40").Value = "Tasa Interna de Retorno TIRM"
$ws.Range("C40").Font.Bold = $true

# Blank "proyectos A/B/C/D" placeholder list (rows 42-45)
$ws.Range("A42").Value = "A"
$ws.Range("A43").Value = "B"
$ws.Range("A44").Value = "C"
$ws.Range("A45").Value = "D"

# --- Worked TIRM example table (rows 47-55) ---
$ws.Range("B47").Value = "Años"
$ws.Range("C47").Value = "A"
$ws.Range("C47").Font.Bold = $true
$ws.Range("D47").Value = "A POS"
$ws.Range("D47").Font.Bold = $true
$ws.Range("E47").Value = "A NEG"
$ws.Range("E47").Font.Bold = $true
$ws.Range("G47").Value = "COC"
$ws.Range("G47").Font.Bold = $true
$ws.Range("H47").Value = 0.1

$ws.Range("B48").Value = 0
$ws.Range("C48").Value = -1500
$ws.Range("E48").Value = -1500

$ws.Range("B49").Value = 1
$ws.Range("C49").Value = 150
$ws.Range("D49").Value = 150

$ws.Range("B50").Value = 2
$ws.Range("C50").Value = 1350
$ws.Range("D50").Value = 1350

$ws.Range("B51").Value = 3
$ws.Range("C51").Value = 150
$ws.Range("D51").Value = 150

$ws.Range("B52").Value = 4
$ws.Range("C52").Value = -150
$ws.Range("E52").Value = -150

$ws.Range("B53").Value = 5
$ws.Range("C53").Value = -600
$ws.Range("E53").Value = -600

# Accounting-style currency format for the cash-flow grid (same as the
# workbook's existing "Currency" cell style).
$ws.Range("C48:E53").NumberFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

# NPV total. Number format set before the formula so the cell lands on the
# built-in "Currency, 2 decimals, red negatives" format (id 8) instead of
# registering a redundant custom format.
$ws.Range("D55").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("D55").Formula = "=NPV(H47,D48:D53)"

# Leave the view scrolled/selected where the author left it.
$ws.Range("E51").Select()
